$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the "Meta description: ..." paragraph that follows the
# "Play Dragon & Phoenix Slot Game for Free - Betsoft" H1 heading.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: insert a new bold paragraph ("Play Dragon & Phoenix Slot Game for
# Free - Betsoft") right before the final paragraph (the DALLE prompt), and
# change that final paragraph's text into the new "Read our review..." text
# while keeping its existing (italic) run formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$startPos = $lastPara.Range.Start

$newHeadingText = "Play Dragon & Phoenix Slot Game for Free - Betsoft"

# Insert the plain text of the new paragraph followed by a paragraph mark.
$insertRange = $d.Range($startPos, $startPos)
$insertRange.InsertAfter($newHeadingText + "`r")

# Bold the text that was just inserted (not the paragraph mark itself).
$boldRange = $d.Range($startPos, $startPos + $newHeadingText.Length)
$boldRange.Font.Bold = 1

# Both the new paragraph and the (now pushed-down) DALLE paragraph need a
# leading empty run to match the document's existing authoring convention,
# so restore it on each using a minimal OOXML fragment insert.
$emptyRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insRangeNew = $d.Range($startPos, $startPos)
$insRangeNew.InsertXML($emptyRunXml)

$newCount = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($newCount)
$dalleStart = $dallePara.Range.Start
$insRangeDalle = $d.Range($dalleStart, $dalleStart)
$insRangeDalle.InsertXML($emptyRunXml)

# Finally, swap the DALLE prompt text for the new review blurb, keeping the
# paragraph's existing (italic) character formatting intact.
$oldText = 'DALLE, please create a feature image for the game "Dragon and Phoenix" that fits the theme and style of the game. The image should be in a cartoon style and should feature a happy Maya warrior with glasses. Make sure the image incorporates elements of the game such as the dragon, the Phoenix bird, the golden money tree, and the Emperor and Empress symbols. The image should be eye-catching and vibrant, and it should make people excited to play the game.'
$newText = 'Read our review of Dragon & Phoenix, the Chinese-themed online slot game from Betsoft. Play for free and enjoy this high-risk, high-reward game.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
